$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 916.5333000000001
$ws.Range("I2").Value = 219
$ws.Range("J2").Value = 1962.8334
$ws.Range("K2").Value = 219
$ws.Range("L2").Value = 1962.8334
$ws.Range("M2").Value = -106
$ws.Range("N2").Value = -2188.8334
$ws.Range("H39").Value = 737.44446
$ws.Range("J39").Value = 999.4
$ws.Range("L39").Value = 2998.2
$ws.Range("N39").Value = -3590.2
$ws.Range("H40").Value = 6025.8237
$ws.Range("I40").Value = 5799.8
$ws.Range("J40").Value = 6348.7144
$ws.Range("K40").Value = 5799.8
$ws.Range("L40").Value = 6348.7144
$ws.Range("M40").Value = -5624.8
$ws.Range("N40").Value = -6698.7144
$ws.Range("H43").Value = 5588.9375
$ws.Range("I43").Value = 1946.25
$ws.Range("J43").Value = 6803.1665
$ws.Range("K43").Value = 1946.25
$ws.Range("L43").Value = 6803.1665
$ws.Range("M43").Value = -1877.25
$ws.Range("N43").Value = -6941.1665
$ws.Range("H70").Value = 6097366
$ws.Range("I70").Value = 105405
$ws.Range("J70").Value = 8702567
$ws.Range("K70").Value = 316215
$ws.Range("L70").Value = 26107701
$ws.Range("M70").Value = -315945
$ws.Range("N70").Value = -26108241
$ws.Range("H73").Value = 6097366
$ws.Range("I73").Value = 105405
$ws.Range("J73").Value = 8702567
$ws.Range("K73").Value = 316215
$ws.Range("L73").Value = 26107701
$ws.Range("M73").Value = -315279
$ws.Range("N73").Value = -26109573
$ws.Range("H86").Value = 4300.4
$ws.Range("I86").Value = 5000.6665
$ws.Range("J86").Value = 3250
$ws.Range("K86").Value = 5000.6665
$ws.Range("L86").Value = 3250
$ws.Range("M86").Value = -3877.6665
$ws.Range("N86").Value = -5496
$ws.Range("H89").Value = 4300.4
$ws.Range("I89").Value = 5000.6665
$ws.Range("J89").Value = 3250
$ws.Range("K89").Value = 25003.3325
$ws.Range("L89").Value = 16250
$ws.Range("M89").Value = -19387.3325
$ws.Range("N89").Value = -27482
$ws.Range("H94").Value = 2718
$ws.Range("I94").Value = 3147.5
$ws.Range("K94").Value = 3147.5
$ws.Range("M94").Value = -2696.5
$ws.Range("H135").Value = 2271.9375
$ws.Range("I135").Value = 2077.4546
$ws.Range("J135").Value = 2699.8
$ws.Range("K135").Value = 18697.0914
$ws.Range("L135").Value = 24298.2
$ws.Range("M135").Value = -16162.0914
$ws.Range("N135").Value = -29368.2
$ws.Range("H141").Value = 2440.4167
$ws.Range("I141").Value = 2148
$ws.Range("K141").Value = 6444
$ws.Range("M141").Value = -1264

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2716.7974
$ws.Range("I32").Value = 2292.3242
$ws.Range("J32").Value = 8999
$ws.Range("K32").Value = 2292.3242
$ws.Range("L32").Value = 8999
$ws.Range("M32").Value = -2005.3242
$ws.Range("N32").Value = -9573
$ws.Range("H74").Value = 30308668
$ws.Range("I74").Value = 66670150
$ws.Range("K74").Value = 66670150
$ws.Range("M74").Value = -66669276
$ws.Range("H77").Value = 30308668
$ws.Range("I77").Value = 66670150
$ws.Range("K77").Value = 333350750
$ws.Range("M77").Value = -333346382

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H102").Value = 26870.334
$ws.Range("I102").Value = 24999.5
$ws.Range("J102").Value = 30612
$ws.Range("K102").Value = 24999.5
$ws.Range("L102").Value = 30612
$ws.Range("M102").Value = -21754.5
$ws.Range("N102").Value = -37102

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H137").Value = 79777
$ws.Range("J137").Value = 79777
$ws.Range("L137").Value = 79777
$ws.Range("N137").Value = -89977
$ws.Range("H138").Value = 62544.168
$ws.Range("J138").Value = 62544.168
$ws.Range("L138").Value = 62544.168
$ws.Range("N138").Value = -72824.16800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H8").Value = 247
$ws.Range("I8").Value = 247
$ws.Range("K8").Value = 741
$ws.Range("M8").Value = -602
$ws.Range("H16").Value = 1888.3334
$ws.Range("I16").Value = 1498.3334
$ws.Range("J16").Value = 2083.3333
$ws.Range("K16").Value = 4495.0002
$ws.Range("L16").Value = 6249.999899999999
$ws.Range("M16").Value = -4322.0002
$ws.Range("N16").Value = -6595.999899999999
$ws.Range("H60").Value = 822414.9
$ws.Range("I60").Value = 180
$ws.Range("J60").Value = 1644649.8
$ws.Range("K60").Value = 540
$ws.Range("L60").Value = 4933949.4
$ws.Range("M60").Value = -289
$ws.Range("N60").Value = -4934451.4
$ws.Range("H94").Value = 8689.286
$ws.Range("J94").Value = 11805.4
$ws.Range("L94").Value = 35416.2
$ws.Range("N94").Value = -36768.2
$ws.Range("H140").Value = 1901.5
$ws.Range("I140").Value = 1231.1428
$ws.Range("K140").Value = 3693.4284
$ws.Range("M140").Value = 1486.5716

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 548.6
$ws.Range("J2").Value = 872.6667
$ws.Range("L2").Value = 872.6667
$ws.Range("N2").Value = -1098.6667
$ws.Range("H80").Value = 2992.5
$ws.Range("I80").Value = 2485
$ws.Range("J80").Value = 3500
$ws.Range("K80").Value = 2485
$ws.Range("L80").Value = 3500
$ws.Range("M80").Value = -1487
$ws.Range("N80").Value = -5496
$ws.Range("H83").Value = 2992.5
$ws.Range("I83").Value = 2485
$ws.Range("J83").Value = 3500
$ws.Range("K83").Value = 12425
$ws.Range("L83").Value = 17500
$ws.Range("M83").Value = -7433
$ws.Range("N83").Value = -27484
$ws.Range("H106").Value = 65624.5
$ws.Range("J106").Value = 65624.5
$ws.Range("L106").Value = 65624.5
$ws.Range("N106").Value = -68148.5
$ws.Range("H109").Value = 70285
$ws.Range("J109").Value = 70285
$ws.Range("L109").Value = 70285
$ws.Range("N109").Value = -72365

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 5000000
$ws.Range("J14").Value = 0
$ws.Range("L14").Value = 0
$ws.Range("N14").ClearContents()
$ws.Range("H135").Value = 63862.11
$ws.Range("J135").Value = 63862.11
$ws.Range("L135").Value = 63862.11
$ws.Range("N135").Value = -74002.11

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 12862.571
$ws.Range("I54").Value = 9999.909
$ws.Range("K54").Value = 9999.909
$ws.Range("M54").Value = -9479.909
$ws.Range("H62").Value = 11129
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 11129
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 11129
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -12377
$ws.Range("H65").Value = 11129
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 11129
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 55645
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -61885
$ws.Range("H104").Value = 9056.333000000001
$ws.Range("J104").Value = 9056.333000000001
$ws.Range("L104").Value = 9056.333000000001
$ws.Range("N104").Value = -16044.333
